$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.464.75"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "1.842.85"
$ws.Range("E3").Value = "  -1.90%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.40"
$ws.Range("E5").Value = "  -5.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5205"
$ws.Range("E7").Value = "  -1.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3278"
$ws.Range("E8").Value = "  -4.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06780"
$ws.Range("E9").Value = "  -2.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.67"
$ws.Range("E10").Value = "  -6.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7734"
$ws.Range("E11").Value = "  -3.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07699"
$ws.Range("E12").Value = "  -0.25%  "

$ws.Range("D13").Value = "1.793.99"
$ws.Range("E13").Value = "  -4.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.25"
$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.030"
$ws.Range("E15").Value = "  -2.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.94"
$ws.Range("E17").Value = "  -4.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9996"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007962"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("D20").Value = "26.444.00"
$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").Value = "2.061.33"
$ws.Range("E21").Value = "  -2.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.579"
$ws.Range("E22").Value = "  -3.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.490"
$ws.Range("E23").Value = "  -5.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.971"
$ws.Range("E24").Value = "  -3.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.25"
$ws.Range("E25").Value = "  -1.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.209"
$ws.Range("E26").Value = "  -7.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.648"
$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.98"
$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.44"
$ws.Range("E29").Value = "  -1.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.203"
$ws.Range("E30").Value = "  -3.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.136"
$ws.Range("E31").Value = "  -4.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08729"
$ws.Range("E32").Value = "  -1.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04806"
$ws.Range("E33").Value = "  -1.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.131"
$ws.Range("E34").Value = "  -3.72%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7094"
$ws.Range("E35").Value = "  -2.48%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.835"
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.086"
$ws.Range("E37").Value = "  -6.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.233"
$ws.Range("E38").Value = "  -4.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01762"
$ws.Range("E39").Value = "  -4.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4850"
$ws.Range("E40").Value = "  -5.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.56"
$ws.Range("E41").Value = "  -4.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8951"
$ws.Range("E42").Value = "  -6.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.066"
$ws.Range("E43").Value = "  -1.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.741"
$ws.Range("E45").Value = "  -4.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4172"
$ws.Range("E46").Value = "  -6.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.003"
$ws.Range("E48").Value = "  -3.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.03"
$ws.Range("E49").Value = "  -3.27%  "

$ws.Range("E50").Value = "  -9.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8887"
$ws.Range("E51").Value = "  +0.66%  "
